$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. "4.2 Items" sheet: fill in previously-blank/zero LastCount values now
#    that the Action column on "4.2 Timestamps" carries explicit counts.
# ---------------------------------------------------------------------------
$itemsWs = $wb.Worksheets.Item("4.2 Items")
$itemsWs.Range("B2").Value = 19
$itemsWs.Range("B8").Value = 6
$itemsWs.Range("C8").Value = 5
$itemsWs.Range("B9").Value = 4
$itemsWs.Range("C9").Value = 3
$itemsWs.Range("B12").Value = 11

# ---------------------------------------------------------------------------
# 2. "4.2 Timestamps" sheet: append the new activity rows (68-89). Each row
#    in the "Action" column now records a count, e.g. "add 1"/"subtract 2",
#    per the commit message.
# ---------------------------------------------------------------------------
$tsWs = $wb.Worksheets.Item("4.2 Timestamps")

$newRows = @(
    @("2024-01-14 14:18:00", "Laptop Charger ",   "add 1",      ""),
    @("2024-01-14 14:18:10", "Laptop Charger ",   "add 5",      ""),
    @("2024-01-14 14:18:26", "Desktop Mini G9",   "add 1",      "SAN124323"),
    @("2024-01-14 14:18:30", "Desktop Mini G9",   "add 1",      "SAN124354"),
    @("2024-01-14 14:18:38", "USB DVD-RW Drive",  "add 2",      ""),
    @("2024-01-14 14:18:48", "USB DVD-RW Drive",  "subtract 2", ""),
    @("2024-01-14 14:22:16", "Desktop Mini G9",   "add 1",      "SAN124354"),
    @("2024-01-14 14:22:43", "Laptop Charger ",   "subtract 1", ""),
    @("2024-01-14 14:28:21", "Desktop Mini G9",   "add",        "SAN123456"),
    @("2024-01-14 14:28:27", "Desktop Mini G9",   "subtract",   "SAN123456"),
    @("2024-01-14 14:47:39", "Desktop Mini G9",   "add 1",      "SAN13579"),
    @("2024-01-14 14:47:48", "Desktop Mini G9",   "add 1",      "SAN13578"),
    @("2024-01-14 14:47:55", "Desktop Mini G9",   "subtract 2", ""),
    @("2024-01-14 14:48:02", "Desktop Mini G9",   "subtract 2", ""),
    @("2024-01-14 14:58:04", "Desktop Mini G9",   "add 1",      "SAN124578"),
    @("2024-01-14 14:58:20", "Desktop Mini G9",   "subtract 1", ""),
    @("2024-01-14 14:58:48", "Desktop Mini G9",   "add 1",      "SAN124589"),
    @("2024-01-14 14:59:05", "Desktop Mini G9",   "subtract 1", ""),
    @("2024-01-14 15:06:14", "Desktop Mini G9",   "add 1",      "SAN223344"),
    @("2024-01-14 15:06:19", "Desktop Mini G9",   "add 1",      "SAN445566"),
    @("2024-01-14 15:06:40", "Desktop Mini G9",   "subtract 1", "SAN223344"),
    @("2024-01-14 15:06:50", "Desktop Mini G9",   "subtract 1", "SAN445566")
)

$row = 68
foreach ($entry in $newRows) {
    $tsWs.Cells.Item($row, 1).Value = $entry[0]
    $tsWs.Cells.Item($row, 2).Value = $entry[1]
    $tsWs.Cells.Item($row, 3).Value = $entry[2]
    if ($entry[3] -ne "") {
        $tsWs.Cells.Item($row, 4).Value = $entry[3]
    }
    $row = $row + 1
}

# ---------------------------------------------------------------------------
# 3. "BR Timestamps" sheet: the trailing blank SAN cell on row 10 was never
#    populated (no SAN for "Laptop Charger") - drop the stray empty cell.
# ---------------------------------------------------------------------------
$brTsWs = $wb.Worksheets.Item("BR Timestamps")
$brTsWs.Cells.Item(10, 4).ClearContents()

# ---------------------------------------------------------------------------
# 4. "All SANs" sheet: remove the duplicate/incorrect SAN106026 entry - the
#    correct "Laptop x360 G8" sighting for that timestamp is SAN106008.
# ---------------------------------------------------------------------------
$sansWs = $wb.Worksheets.Item("All SANs")
$sansWs.Rows.Item(55).Delete()
